# Append a trailing "." to the "Answer = 3 + 2 + 1" line in the
# Subtitle placeholder of slide 1 (ProblemStatements/FindPairsCarryForward).
#
# Commit message: "Checking commit and push from laptop." — a small,
# incidental text tweak made right before a resave (the rest of the
# upstream diff is just PowerPoint's own round-trip noise: dropped
# cached field text, stripped extLst/creationId metadata, etc., which
# are not reachable/meaningful through the PowerPoint object model).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item(2)            # "Subtitle 2" placeholder
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(11)            # "Answer = 3 + 2 + 1"

if ($para.Text.TrimEnd() -eq "Answer = 3 + 2 + 1") {
    [void]$para.InsertAfter(".")
}
